# Re-colour the deck's theme: swap the custom "Integral" palette for the
# built-in "Office Theme" palette (the one already used by the Notes
# Master), by rewriting every slot of the Slide Master's theme colour
# scheme to the standard Office RGB values.
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Office theme palette, in ThemeColorScheme index order:
#  1 dk1      000000
#  2 lt1      FFFFFF
#  3 dk2      44546A
#  4 lt2      E7E6E6
#  5 accent1  5B9BD5
#  6 accent2  ED7D31
#  7 accent3  A5A5A5
#  8 accent4  FFC000
#  9 accent5  4472C4
# 10 accent6  70AD47
# 11 hlink    0563C1
# 12 folHlink 954F72
$officeRGB = @(
    0,
    16777215,
    6968388,
    15132391,
    13998939,
    3243501,
    10855845,
    49407,
    12874308,
    4697456,
    12673797,
    7491477
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeRGB[$i - 1]
}
